# --- Adding "Acceptance test time report" sheet (acceptance tests for the time report feature) ---
$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Put the selection on sheet 2 where it ends up after the edit, before we move away from it
$ws2.Activate()
$ws2.Range("E5").Select()

# Insert the new worksheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Acceptance test time report"

# Match the sheet tab color used for this new acceptance-test sheet
$ws3.Tab.Color = 5287936

# Header row - reuse the same header styling as the other acceptance test sheet
$ws2.Range("A1:E1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)
$ws3.Range("A1").Value = "Test ID"
$ws3.Range("B1").Value = "Description"
$ws3.Range("C1").Value = "Expected Result"
$ws3.Range("D1").Value = "Precondition"
$ws3.Range("E1").Value = "Comments"

# Test case rows (content entered in the same order as originally authored,
# so that rows sharing identical text reuse the same shared-string entry)
$ws3.Range("B2").Value = "the maneger want to see the time report `nhe enter to  GUI and `nenter `nvalid datetime to start the report`nvalid datetime to end the report`nend click enter "
$ws3.Range("C2").Value = "the system enter the the database and get all the informatin of the Orders and queue form the db `nthe system analize the data.`nSplit it to how much was make orders, and how mush was enter the queue`nhow much from the orders was late `nthe system make report of it with grafics and sent it as a result to the manager "
$ws3.Range("D2").Value = "the server is accisble `nthe DB is accesable `nthe data of the Orders in that time range aviable`nthe data of the Queue in that time range aviable"
$ws3.Range("E2").Value = "the manager want to see report about the orders time and about the queue times`nthe system get to the report "
$ws3.Range("A2").Value = "reportWithDataSusssesful"
$ws3.Range("A3").Value = "reportWithoutDataFail"
$ws3.Range("C3").Value = "the system enter the the database and get all the informatin of the Orders and queue form the db `nthe system find that there no data in that time range and sent messege : `"data not found for this range`""
$ws3.Range("D3").Value = "the server is accisble `nthe DB is accesable `n"
$ws3.Range("E3").Value = "the manager want to see report about the orders time and about the queue times`nthe system cant find data on this time range`nthe system get to the report "
$ws3.Range("A4").Value = "reportWithoutDateTimeFail"
$ws3.Range("B4").Value = "the maneger want to see the time report `nhe enter to  GUI and `nenter `nvalid datetime to end the report`nend click enter "
$ws3.Range("C4").Value = "the GUI check if all the fuilds are full. The system write massege that need to fill the datetime to start report"
$ws3.Range("E4").Value = "the manager want to see report about, but forget to fill all the fields so ther GUI will tell them that"

# Alignment: Test ID column is centered, the rest are left-aligned with wrap
$ws3.Range("A2").HorizontalAlignment = -4108
$ws3.Range("A2").VerticalAlignment = -4108
$ws3.Range("B2:E2").HorizontalAlignment = -4131
$ws3.Range("B2:E2").VerticalAlignment = -4108
$ws3.Range("B2:E2").WrapText = $true
$ws3.Range("A3").HorizontalAlignment = -4108
$ws3.Range("A3").VerticalAlignment = -4108
$ws3.Range("B3:E3").HorizontalAlignment = -4131
$ws3.Range("B3:E3").VerticalAlignment = -4108
$ws3.Range("B3:E3").WrapText = $true
$ws3.Range("A4").HorizontalAlignment = -4108
$ws3.Range("A4").VerticalAlignment = -4108
$ws3.Range("B4:E4").HorizontalAlignment = -4131
$ws3.Range("B4:E4").VerticalAlignment = -4108
$ws3.Range("B4:E4").WrapText = $true

# Row heights (matches the authored sheet)
$ws3.Rows.Item(1).RowHeight = 17
$ws3.Rows.Item(2).RowHeight = 136
$ws3.Rows.Item(3).RowHeight = 119
$ws3.Rows.Item(4).RowHeight = 102

# Column widths (matches the authored sheet)
$ws3.Columns.Item(1).ColumnWidth = 35.25
$ws3.Columns.Item(2).ColumnWidth = 31.42
$ws3.Columns.Item(3).ColumnWidth = 49.09
$ws3.Columns.Item(4).ColumnWidth = 27.25
$ws3.Columns.Item(5).ColumnWidth = 30.75

# View settings for the new sheet
$ws3.Activate()
$ws3.Application.ActiveWindow.Zoom = 194
$ws3.Range("C14").Select()

